$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map the existing numeric "RESPUESTA" values (1,2,3,4) in column F
# to their letter equivalents (a,b,c,d) as shared-string text values.
$letters = @{ 1 = "a"; 2 = "b"; 3 = "c"; 4 = "d" }

for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $num = [int]$cell.Value2
    $cell.Value = $letters[$num]
}

# Add a new empty row below the table (F17) with an underlined font style,
# matching the new cellXfs entry / font added for this cell.
$newCell = $ws.Cells.Item(17, 6)
$newCell.Font.Underline = $true

# Move the active selection to the newly added cell.
$ws.Range("F17").Select()
